# ---------------------------------------------------------------------------
# Design document edit: relocate the "_GoBack" bookmark.
#
# The paragraph's visible text does not change at all -- only where the
# (invisible, zero-length) "_GoBack" bookmark sits inside it changes, which
# in turn changes where Word has to split the underlying runs:
#
#   BEFORE: "...communicate with other subsystems in that <i>partition</i>...
#            ...edit a recipe that doesn't belong to them). Thus, the
#            following reasons why we went with a combination[BOOKMARK] of
#            both layering and partitioning when choosing a system design."
#
#   AFTER:  "...communicate with other[BOOKMARK] subsystems in that
#            <i>partition</i>...edit a recipe that doesn't belong to them).
#            Thus, the following reasons why we went with a combination of
#            both layering and partitioning when choosing a system design."
#
# So two things must happen:
#   1) The run that used to read "...communicate with other subsystems in
#      that " must be split so the bookmark sits right after "...other".
#   2) The two runs that used to be split by the bookmark's old location
#      ("...a combination" / " of both...design.") must become one run,
#      since nothing separates them any more.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the relevant anchor positions dynamically via Find rather than
# hard-coding character offsets.
# ---------------------------------------------------------------------------

# (A) New split point for the bookmark: right after "...communicate with
#     other" / before " subsystems in that ...".
$findA = $d.Content
$findA.Find.ClearFormatting()
$okA = $findA.Find.Execute("communicate with other", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okA) { throw "Could not find 'communicate with other'" }
$splitPos = $findA.End

# (B) Boundary that must NOT be disturbed: right after "...belong to them)."
#     / before " Thus, the following...". This sits just before the two
#     runs that need to merge, and must stay separate from them.
$findB = $d.Content
$findB.Find.ClearFormatting()
$okB = $findB.Find.Execute("belong to them).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okB) { throw "Could not find 'belong to them).'" }
$barrierPos = $findB.End

# ---------------------------------------------------------------------------
# Part 1 -- move the "_GoBack" bookmark to its new location (A). Bookmarks
# can't have their Start/End reassigned directly, so delete + re-add.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$newBmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $newBmRange)

# Inserting the bookmark mid-run splits its <w:t> in two; the leading half
# keeps a stray xml:space="preserve" it no longer needs. Force a narrow,
# scoped self Find & Replace over just that shortened run so it gets
# re-serialized cleanly (the range is kept tight so nothing else nearby is
# touched).
$cleanupRange = $d.Range($splitPos - 30, $splitPos + 5)
$cleanupRange.Find.ClearFormatting()
$null = $cleanupRange.Find.Execute("communicate with other", $true, $false, $false, $false, $false, $true, 1, $false, "communicate with other", 2)

# ---------------------------------------------------------------------------
# Part 2 -- merge the two runs that used to be separated by the bookmark's
# old location ("...a combination" + " of both...design."). A temporary
# bookmark is planted at boundary (B) first, acting as a hard barrier so the
# forced rewrite below only fuses those two runs and does not also swallow
# the preceding, identically-formatted run.
# ---------------------------------------------------------------------------
$barrierRange = $d.Range($barrierPos, $barrierPos)
$d.Bookmarks.Add("TempBarrier", $barrierRange)

$mergeRange = $d.Range($barrierPos, $barrierPos + 60)
$mergeRange.Find.ClearFormatting()
$null = $mergeRange.Find.Execute("combination of both", $true, $false, $false, $false, $false, $true, 1, $false, "combination of both", 2)

$d.Bookmarks.Item("TempBarrier").Delete()
